$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)

$ws.Range("A252").Value = 44326
$ws.Range("B252").Value = 3
$ws.Range("C252").Value = 18
$ws.Range("D252").Value = 182.5187588724397

$ws.Range("A253").Value = 44327
$ws.Range("B253").Value = 1
$ws.Range("C253").Value = 19
$ws.Range("D253").Value = 192.6586899209085

$ws.Range("A254").Value = 44328
$ws.Range("B254").Value = 1
$ws.Range("C254").Value = 19
$ws.Range("D254").Value = 192.6586899209085

$ws.Range("A255").Value = 44329
$ws.Range("B255").Value = 1
$ws.Range("C255").Value = 16
$ws.Range("D255").Value = 162.2388967755019
